# Adds a new "22. 6. 2021" data column to both sheets ("data" and "pocetR")
# and updates the two "aktualizace" footer notes from 1. 6. 2021 to 28. 6. 2021.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data" (sheet1): new column AF, rows 1-75 (header + 74 data rows)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Header cell AF1 - copy style from the existing last header cell (AE1)
$ws1.Range("AE1").Copy()
$ws1.Range("AF1").PasteSpecial(-4122)
$ws1.Range("AF1").Value = "22. 6. 2021"

$data1 = @(0.32,0.31,0.08,0.33,0.21,0.43,0.58,0.19,0.31,0.26,0.29,0.34,0.4,0.37,0.13,0.33,0.14,0.4,0.24,0.47,0.49,0.2,0.67,0.19,0.74,0.15,0.22,0.25,0.31,0.31,0.4,0.32,0.3,0.38,0.34,0.33,0.3,0.29,0.11,0.36,0.2,0.47,0.66,0.16,0.04,0.29,0.22,0.38,0.51,0.21,0.31,0.31,0.29,0.32,0.32,0.25,0.43,0.34,0.43,0.2,0.49,0.3,0.31,0.34,0.39,0.32,0.21,0.36,0.4,0.3,0.34,0.33,0.29,0.31)

$r = 2
foreach ($v in $data1) {
    $ws1.Cells.Item($r, 32).Value = $v
    $r = $r + 1
}

# Footer note (row 76, column A)
$ws1.Range("A76").Value = "Život během pandemie, Imunizace, % respondentů celkově a ve skupinách, aktualizace 28. 6. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR" (sheet2): new column AE, rows 1-38 (header + 37 data rows)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AE1 - copy style from the existing last header cell (AD1)
$ws2.Range("AD1").Copy()
$ws2.Range("AE1").PasteSpecial(-4122)
$ws2.Range("AE1").Value = "22. 6. 2021"

$data2 = @(1904,459,699,746,867,661,376,433,468,258,434,193,118,194,753,580,249,937,967,239,365,334,220,335,411,983,447,223,251,304,237,305,299,544,377,420,1107)

$r = 2
foreach ($v in $data2) {
    $ws2.Cells.Item($r, 31).Value = $v
    $r = $r + 1
}

# Footer note (row 39, column A)
$ws2.Range("A39").Value = "Život během pandemie, Imunizace, velikost dotázaného souboru celkově a ve skupinách, aktualizace 28. 6. 2021"
